# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Cell -> new value mapping (row => new F value)
$updates = @{
    2 = 1391
    3 = 2314
    4 = 402
    6 = 6470
    7 = 323
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
